# Reliability Test 50 Trials 1mL.xlsx
# "added 50 trials for commercial pipette"
#
# - Rename Sheet1 -> "Test Pipette"
# - Add a new sheet "Commericial Pipette" (after "Test Pipette") with the
#   same A/C/D layout (Trial #, dispensed amount, NORMDIST) populated with
#   the 50 commercial-pipette trial readings, plus Mean / StDev / Std Error
#   summary rows 54-56.
# - Leave the new sheet active/selected (mirrors the author tabbing over to
#   review their new data entry) and restore a "select-all-data" selection
#   on the original sheet.

$wb = $excel.ActiveWorkbook

$wsTest = $wb.Worksheets.Item(1)
$wsTest.Name = "Test Pipette"

# New sheet goes right after "Test Pipette".
$wsCommercial = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTest)
$wsCommercial.Name = "Commericial Pipette"

# Column widths matching the source sheet's B/C/D columns (closest
# achievable input given the host's column-width quantization).
$wsCommercial.Columns.Item(2).ColumnWidth = 25.95
$wsCommercial.Columns.Item(3).ColumnWidth = 51.3
$wsCommercial.Columns.Item(4).ColumnWidth = 24.15

# Headers (reuse the same text as "Test Pipette" - engine will dedupe the
# shared-string table automatically).
$wsCommercial.Range("A1").Value = "Trial"
$wsCommercial.Range("C1").Value = "Amount Dispensed Experimental (1mL)Pipette (g)"
$wsCommercial.Range("D1").Value = "Normal Distribution"

# 50 trials of commercial-pipette dispensed amounts (grams).
$amounts = @(0.2024,0.20469999999999999,0.2039,0.2046,0.20369999999999999,0.2044,0.20469999999999999,0.2044,0.20300000000000001,0.2031,0.20399999999999999,0.2041,0.2049,0.2041,0.2039,0.20419999999999999,0.20430000000000001,0.20380000000000001,0.20430000000000001,0.20380000000000001,0.2039,0.20380000000000001,0.2036,0.2041,0.20319999999999999,0.20269999999999999,0.2039,0.20399999999999999,0.2031,0.20419999999999999,0.2039,0.20380000000000001,0.2036,0.2046,0.2034,0.20319999999999999,0.20369999999999999,0.2034,0.2034,0.20419999999999999,0.2041,0.2039,0.20380000000000001,0.2034,0.20449999999999999,0.20330000000000001,0.2039,0.20369999999999999,0.2041,0.2039)

for ($i = 0; $i -lt $amounts.Count; $i++) {
    $row = $i + 2
    $wsCommercial.Cells.Item($row, 1).Value = $i + 1
    $wsCommercial.Cells.Item($row, 3).Value = $amounts[$i]
    $wsCommercial.Cells.Item($row, 4).Formula = "=NORMDIST(C:C,C54,C55,TRUE )"
}

# Summary rows: Mean, Standard Deviation, Standard Error.
$wsCommercial.Range("A54").Value = "Mean "
$wsCommercial.Range("C54").Formula = "=AVERAGE(C2:C51)"

$wsCommercial.Range("A55").Value = "Standard Deviation"
$wsCommercial.Range("C55").Formula = "=STDEV(C2:C51)"

$wsCommercial.Range("A56").Value = "Standard Error"
$wsCommercial.Range("C56").Formula = "=(C55/(SQRT(50)))"

# Selections: new sheet ends up active with C51 selected (last entered
# data cell); original sheet keeps a "select all data" range.
$wsTest.Range("A1:D56").Select() | Out-Null
$wsCommercial.Range("C51").Select() | Out-Null
$wsCommercial.Activate()
